$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 68

$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($r, 3).Value = "Bíobío"
$ws.Cells.Item($r, 4).Value = 45267
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(67, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = 8
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100107
$ws.Cells.Item($r, 8).Value = "Otros"
$ws.Cells.Item($r, 9).Value = 100107002
$ws.Cells.Item($r, 10).Value = "Chirimoya"
$ws.Cells.Item($r, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 19000
$ws.Cells.Item($r, 15).Value = 19000
$ws.Cells.Item($r, 16).Value = 19000
$ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($r, 19).Value = 1900
$ws.Cells.Item($r, 20).Value = 10
